$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove C2 and E2 entirely
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: remove C3 entirely; update E3
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = 4.26438011980097

# Row 4: update C4
$ws.Range("C4").Value = -0.8792832172735632

# Row 5: update C5, E5
$ws.Range("C5").Value = 0.9337833426867226
$ws.Range("E5").Value = 2.730731696345212

# Row 6: update C6, E6
$ws.Range("C6").Value = 2.791140000794257
$ws.Range("E6").Value = 1.68385714213084

# Row 7: update C7
$ws.Range("C7").Value = 0.4451370000809529

# Row 8: update E8
$ws.Range("E8").Value = 2.553470871380514

# Row 9: update E9
$ws.Range("E9").Value = -0.4341460075841019

# Row 10: update E10
$ws.Range("E10").Value = 2.383242923544548

# Row 11: update C11
$ws.Range("C11").Value = 2.2044495746113

# Row 12: update E12
$ws.Range("E12").Value = 2.066615940231942

# Row 13: update E13
$ws.Range("E13").Value = 3.086275812215322

# Row 14: update E14
$ws.Range("E14").Value = -1.135072001636317

# Row 15: update C15, E15
$ws.Range("C15").Value = -3.258619210312896
$ws.Range("E15").Value = 11.44905912635792

# Row 16: update C16
$ws.Range("C16").Value = 0.4255262881966759

# Row 18: update C18, E18
$ws.Range("C18").Value = -0.2814561130375703
$ws.Range("E18").Value = -1.104428907745314

# Row 19: update C19, E19
$ws.Range("C19").Value = -0.6470065423293758
$ws.Range("E19").Value = 3.449881734069282
